# Refresh the legacy GSC "Coverage" export: the rolling date window advanced
# by one day. The oldest day (2025-10-21, row 2) drops off the top, every
# remaining row shifts up by one, and the two newest days (now rows 3 and 4)
# have not been crawled yet so their "Not indexed"/"Indexed" counts are
# blank until the next export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the oldest day's row; Excel shifts rows 3..90 up into 2..89 and
# shrinks the sheet's used range (A1:D90 -> A1:D89) automatically.
$ws.Rows.Item(2).Delete()

# The two newest days (now rows 3 and 4 - "2025-10-23" and "2025-10-24")
# don't have coverage data yet, same as the newest day (row 2, "2025-10-22").
$ws.Range("B2:C4").Value = ""
